# Added ability to ref note in table headings and index
# - "All" heading gets a footnote reference (3)
# - "Sepal Length" heading gets a footnote reference (4)
# - Two new footnotes are appended to the notes/index block at the
#   bottom of the table, pushing the final "no reference" note down
#   to make room for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "All" -> "All(3)" (table column header spanning C6:F6, stored value on C6)
$ws.Range("C6").Value = "All(3)"

# 2. "Sepal Length" -> "Sepal Length(4)" (row group label, repeated on the
#    "Mean" and "Median" rows of that group: A11 and A12)
$ws.Range("A11").Value = "Sepal Length(4)"
$ws.Range("A12").Value = "Sepal Length(4)"

# 3. The last row of the notes/index (A19) becomes the note explaining
#    reference (3)
$ws.Range("A19").Value = "(3: All species of the Iris genus)"

# 4. Insert a new note row (20) explaining reference (4), matching the
#    footnote formatting already used by the notes block (copy from A19)
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "(4: Length of the largest sepal)"

# 5. Insert a new row (21) carrying the general note that used to live in
#    row 19, now pushed to the end of the notes/index block
$ws.Range("A19").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "(This note hath no reference.)"

# Cosmetic: the workbook was re-saved with gridlines hidden and the
# selection reset to the top-left cell
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A1").Select() | Out-Null

Write-Output "Edit applied"
